# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.222.29"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "2.952.07"
$ws.Range("E3").Value = "  +1.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("E7").Value = "  -0.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +0.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.25%  "

$ws.Range("E11").Value = "  +0.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0840"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.47%  "

$ws.Range("D13").Value = "3.418.40"
$ws.Range("E13").Value = "  +1.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.29%  "

$ws.Range("D16").Value = "2.944.36"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.958"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.16%  "

$ws.Range("D18").Value = "51.212.75"
$ws.Range("E18").Value = "  +0.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.96%  "

$ws.Range("E22").Value = "  +1.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.90%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("E25").Value = "  +3.75%  "

$ws.Range("E26").Value = "  -1.81%  "

$ws.Range("E27").Value = "  +1.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +17.05%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  -0.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.111"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.59%  "

$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.40%  "

$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.33%  "

$ws.Range("E36").Value = "  +5.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("E38").Value = "  -1.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.25%  "

$ws.Range("E40").Value = "  -4.05%  "

$ws.Range("E41").Value = "  -0.97%  "

$ws.Range("E42").Value = "  +1.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "123.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.74%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.39%  "

$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.278"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +15.65%  "

$ws.Range("E47").Value = "  +0.76%  "

$ws.Range("D48").Value = "2.033.26"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("E49").Value = "  +0.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0347"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.93%  "
